# Update cryptos list (Price / Volume(1h) columns) with freshly scraped
# values, as produced by the scheduled GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price cells are plain text ("28.903.32", "1.000", ...), so stage
# them through a Text number format while writing -- this stops Excel
# from coercing look-alike numbers (e.g. "1.000" -> 1) -- then restore
# the default "Normal" style so the cell keeps its original (unstyled)
# formatting, matching the rest of the sheet.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.903.32"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.16%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.831.99"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.53%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9992"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.04%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6937"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.68%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.000"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07669"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.99%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3040"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.64%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.25"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.36%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07808"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.12%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "93.01"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.11%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.830.15"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.76%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.093"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.78%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6814"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.80%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.505"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.16%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008238"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.34%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "28.906.51"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.23%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "240.23"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.22%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.074.00"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.97%  "

$ws.Range("E21").Value = "  -2.14%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.03%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.447"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.71%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.001"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.12%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1496"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.77%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "158.32"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.61%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.731"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.91%  "

$ws.Range("E28").Value = "  -2.32%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.538"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.45%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.226"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.37%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.133"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.46%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.192"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.33%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05113"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.30%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7732"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.50%  "

$ws.Range("E35").Value = "  -1.69%  "

$ws.Range("E36").Value = "  -3.29%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.696"
$ws.Range("D37").Style = "Normal"

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.276.96"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.80%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01856"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.12%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.694"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.71%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9551"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.95%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.141"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.00%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "106.97"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.24%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.000"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.05%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "9.677"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.20%  "

$ws.Range("E46").Value = "  -1.22%  "

$ws.Range("E47").Value = "  -0.24%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.974.82"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.87%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "63.69"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -7.68%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.749"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.27%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.962"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.77%  "
